$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (Text number format) on D-column cells whose new value would
# otherwise be auto-parsed by Excel as a number (loses the original text
# formatting, e.g. trailing zeros like "225.10" or leading-zero decimals).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Update Price (D) and Volume/1h (E) columns with the latest scraped values.
$ws.Range("D2").Value = "34.338.55"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "1.796.89"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "225.10"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "0.598"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "36.29"
$ws.Range("E8").Value = "  +3.55%  "
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").Value = "0.0675"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").Value = "0.0963"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").Value = "2.054.46"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "11.28"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "1.806.80"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "0.633"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "34.295.15"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "4.42"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "68.58"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("D19").Value = "244.16"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "0.0₃0771"
$ws.Range("E20").Value = "  -3.75%  "
$ws.Range("D21").Value = "11.35"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "4.07"
$ws.Range("E23").Value = "  -2.74%  "
$ws.Range("D24").Value = "2.19"
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").Value = "170.76"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("D26").Value = "7.90"
$ws.Range("E26").Value = "  +4.83%  "
$ws.Range("D27").Value = "17.29"
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "1.23"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").Value = "3.78"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "3.88"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "0.0512"
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("D34").Value = "1.77"
$ws.Range("E34").Value = "  -4.19%  "
$ws.Range("D35").Value = "1.360.05"
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("D36").Value = "0.645"
$ws.Range("E36").Value = "  -5.07%  "
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  -7.98%  "
$ws.Range("D39").Value = "0.0185"
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").Value = "80.63"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").Value = "0.936"
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("E44").Value = "  +4.70%  "
$ws.Range("D45").Value = "13.17"
$ws.Range("E45").Value = "  -5.19%  "
$ws.Range("D46").Value = "0.0496"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("D47").Value = "1.956.15"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("E48").Value = "  -5.39%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "101.86"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").Value = "0.0₆0121"
$ws.Range("E51").Value = "  -7.32%  "
